$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45067.52, 6, 2025, "06/2025"),
    @(3, 35051.88, 6, 2025, "06/2025"),
    @(4, 11320.2, 6, 2025, "06/2025"),
    @(5, 12749.74, 6, 2025, "06/2025"),
    @(6, 12209.32, 6, 2025, "06/2025"),
    @(9, 10140.2, 6, 2025, "06/2025"),
    @(2, 20463.44, 5, 2025, "05/2025"),
    @(5, 29720.49, 5, 2025, "05/2025"),
    @(6, 19666.55, 5, 2025, "05/2025"),
    @(7, 21883.41, 5, 2025, "05/2025"),
    @(8, 45418.95, 5, 2025, "05/2025"),
    @(9, 39575.07, 5, 2025, "05/2025"),
    @(12, 15499.7, 5, 2025, "05/2025"),
    @(13, 10155.95, 5, 2025, "05/2025"),
    @(14, 34413.69, 5, 2025, "05/2025"),
    @(15, 28325.85, 5, 2025, "05/2025"),
    @(16, 12000.74, 5, 2025, "05/2025"),
    @(19, 17756.75, 5, 2025, "05/2025"),
    @(20, 23593.88, 5, 2025, "05/2025"),
    @(21, 9475.47, 5, 2025, "05/2025"),
    @(22, 27766.95, 5, 2025, "05/2025"),
    @(23, 11639.15, 5, 2025, "05/2025"),
    @(26, 28728.43, 5, 2025, "05/2025"),
    @(27, 15807.63, 5, 2025, "05/2025"),
    @(28, 24406.06, 5, 2025, "05/2025"),
    @(29, 30466.72, 5, 2025, "05/2025"),
    @(30, 23720.49, 5, 2025, "05/2025"),
    @(1, 31020.03, 4, 2025, "04/2025"),
    @(2, 48732.41, 4, 2025, "04/2025"),
    @(3, 35050.25, 4, 2025, "04/2025"),
    @(4, 7210.62, 4, 2025, "04/2025"),
    @(7, 10042.97, 4, 2025, "04/2025"),
    @(8, 7630.52, 4, 2025, "04/2025"),
    @(9, 19095.69, 4, 2025, "04/2025"),
    @(10, 16812.5, 4, 2025, "04/2025"),
    @(11, 6701.92, 4, 2025, "04/2025"),
    @(14, 38115.7, 4, 2025, "04/2025"),
    @(15, 12545.33, 4, 2025, "04/2025"),
    @(16, 14981.31, 4, 2025, "04/2025"),
    @(17, 14359.25, 4, 2025, "04/2025"),
    @(22, 18901.75, 4, 2025, "04/2025"),
    @(23, 22062.13, 4, 2025, "04/2025"),
    @(24, 17802.7, 4, 2025, "04/2025"),
    @(25, 6150.450000000001, 4, 2025, "04/2025"),
    @(28, 31344.88, 4, 2025, "04/2025"),
    @(29, 33546.1, 4, 2025, "04/2025"),
    @(30, 32403.27, 4, 2025, "04/2025"),
    @(3, 21522.57, 3, 2025, "03/2025"),
    @(6, 37624.31, 3, 2025, "03/2025"),
    @(7, 13035.05, 3, 2025, "03/2025"),
    @(10, 20789.35, 3, 2025, "03/2025"),
    @(11, 10330.44, 3, 2025, "03/2025"),
    @(12, 25429.94, 3, 2025, "03/2025"),
    @(13, 14198.2, 3, 2025, "03/2025"),
    @(14, 21391.35, 3, 2025, "03/2025"),
    @(17, 26250.86, 3, 2025, "03/2025"),
    @(18, 3373.67, 3, 2025, "03/2025"),
    @(19, 17408.6, 3, 2025, "03/2025"),
    @(20, 5005.650000000001, 3, 2025, "03/2025"),
    @(21, 22003.87, 3, 2025, "03/2025"),
    @(24, 30411.66, 3, 2025, "03/2025"),
    @(25, 27184.21, 3, 2025, "03/2025"),
    @(26, 19290.09, 3, 2025, "03/2025"),
    @(27, 13162.5, 3, 2025, "03/2025"),
    @(28, 11137.4, 3, 2025, "03/2025"),
    @(31, 41403.64, 3, 2025, "03/2025")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

Write-Host "done"
